# Method for detecting missing files
# Applies:
#  - sharedStrings text fixes (Sheet2/"stm_sxm" & two CV labels)
#  - two new rows (11 & 12) on Sheet3 listing test_missing_1 / test_missing_2
#  - updated view/selection state on both sheets

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- Text corrections (shared strings) -----------------------------------
$ws2.Range("B16").Value = "stm-nanonis-sxm"
$ws2.Range("B21").Value = "CV_153505_ 1"
$ws2.Range("B22").Value = "CV_153605_ 2"

# --- New rows describing the missing-file detection method ---------------
# A11 stays empty but keeps the sheet's default (non wrap-text) style.
$a11 = $ws3.Range("A11")
$a11.WrapText = $false
$a11.VerticalAlignment = -4107
$a11.Font.Name = "Calibri"
$a11.Font.Size = 11

$ws3.Range("B11").Value = "test_missing_1"
$ws3.Range("B12").Value = "test_missing_2"

# --- View / selection state ------------------------------------------------
# Select Sheet3's cell first, then finish on Sheet2 so it remains the
# active/selected tab (matches the original file where Sheet2 is tabSelected).
$ws3.Activate()
$ws3.Range("B14").Select()

$ws2.Activate()
$ws2.Range("C22").Select()
